# Applies the Nov-2023 performance-sheet touch-up to "绩效表":
#   - K8 / K9 / K11: replace the #N/A error result with a blank (empty) cell
#   - Row 12: fill in a new line item for 黄礼闯 (曹卓补充订单 / 其他业务, 待完成)
#   - F30 / I30: bump the "软性指标合计" tally from 3 to 4 inputs

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("绩效表")

# --- Clear the stray #N/A results in the 绩效系数 column ---
$ws.Range("K8").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("K11").Value = ""

# --- Row 12: new entry for 黄礼闯 ---
$ws.Range("A12").Value = "黄礼闯"
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = "曹卓补充订单"
$ws.Range("D12").Value = "其他业务"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = "曹卓交付三个订单所需数据"
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = "待完成"
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""

# --- 软性指标合计 row (30): one more NA term folded into the tally ---
$ws.Range("F30").Value = 4
$ws.Range("I30").Value = "NA+NA+NA+NA=NA"
